$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (KB / 2024-04-18 / 제일엠앤에스 entry) - all subsequent rows shift up by one.
$ws.Rows.Item(7).Delete()

# After the shift, the former rows 23 (한국제14호스팩) and 24 (디앤디파마텍) are now
# rows 22 and 23 respectively, but the target layout has them swapped back
# (디앤디파마텍 on row 22, 한국제14호스팩 on row 23). Write the swapped rows explicitly
# (date-like text cells are apostrophe-prefixed so they stay text, not real dates).
$ws.Range("A22").Value = "한국"
$ws.Range("B22").Value = "'2024-04-22"
$ws.Range("C22").Value = "디앤디파마텍"
$ws.Range("D22").Value = "한국"
$ws.Range("E22").Value = "한국"
$ws.Range("F22").Value = "'2024-04-25"
$ws.Range("G22").Value = "'2024-05-02"
$ws.Range("H22").Value = 36300
$ws.Range("I22").Value = 1100000
$ws.Range("J22").Value = 33000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100

$ws.Range("A23").Value = "한국"
$ws.Range("B23").Value = "'2024-06-10"
$ws.Range("C23").Value = "한국제14호스팩"
$ws.Range("D23").Value = "한국"
$ws.Range("E23").Value = "한국"
$ws.Range("F23").Value = "'2024-06-13"
$ws.Range("G23").Value = "'2024-06-19"
$ws.Range("H23").Value = 8000
$ws.Range("I23").Value = 4000000
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100
